$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2022 column (column R), mirroring column Q's formatting
# (copy/paste-formats keeps the existing cellXfs entries instead of
# minting new, duplicate style records).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)

$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)

$ws.Range("R4").Value = 2022
$ws.Range("R5").Value = 8.6821914120339212
$ws.Range("R6").Value = 12.221423436376707

# Update selection to match the recorded active cell after the edit.
$ws.Range("S4").Select()
